# ERP-435 - Separate the Manchester/Glasgow tribunal addresses onto multiple
# lines (address line 1/2/3, town, postcode) instead of one long comma
# separated address string, so that generated letters show the address
# across multiple lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop any existing hyperlinks up front - the email hyperlinks currently
# anchored at B7 / B12 will be re-created at their new locations (B11 / B19)
# once the rows below have been re-laid-out. (Hyperlinks.Delete() clears all
# hyperlinks on the sheet regardless of which range it is invoked from.)
$ws.Range("A1").Hyperlinks.Delete() | Out-Null

# Field name / value pairs for rows 3 through 19 (rows 1 & 2 - the FIELDS /
# VALUES header and positionType row - are unchanged). Writing them in this
# top-to-bottom, column-A-then-B order keeps new shared-string creation
# order aligned with row order.
$fieldRows = @(
  @{ Field = "tribunalManchesterAddressLine1"; Value = "Manchester Employment Tribunal," },
  @{ Field = "tribunalManchesterAddressLine2"; Value = "Alexandra House," },
  @{ Field = "tribunalManchesterAddressLine3"; Value = "14-22 The Parsonage," },
  @{ Field = "tribunalManchesterTown";         Value = "Manchester," },
  @{ Field = "tribunalManchesterPostCode";     Value = "M3 2JA" },
  @{ Field = "tribunalManchesterTelephone";    Value = "0161 833 6100" },
  @{ Field = "tribunalManchesterFax";          Value = "0870 739 4433" },
  @{ Field = "tribunalManchesterDX";           Value = "DX 743570" },
  @{ Field = "tribunalManchesterEmail";        Value = "Manchesteret@justice.gov.uk" },
  @{ Field = "tribunalGlasgowAddressLine1";    Value = "Eagle Building," },
  @{ Field = "tribunalGlasgowAddressLine2";    Value = "215 Bothwell Street," },
  @{ Field = "tribunalGlasgowTown";            Value = "Glasgow," },
  @{ Field = "tribunalGlasgowPostCode";        Value = "G2 7TS" },
  @{ Field = "tribunalGlasgowTelephone";       Value = "0141 204 0730" },
  @{ Field = "tribunalGlasgowFax";             Value = "01264 785 177" },
  @{ Field = "tribunalGlasgowDX";              Value = "DX 580003" },
  @{ Field = "tribunalGlasgowEmail";           Value = "glasgowet@justice.gov.uk" }
)

$rowNum = 3
foreach ($entry in $fieldRows) {
  $ws.Cells.Item($rowNum, 1).Value = $entry.Field
  $ws.Cells.Item($rowNum, 2).Value = $entry.Value
  $rowNum = $rowNum + 1
}

# The field-name column wraps onto multiple lines for the Manchester address
# block (rows 2-6: positionType plus the three new address lines and town),
# matching the rest of the FIELDS column formatting.
$ws.Range("A2:A6").WrapText = $true

# Keep the trailing right-aligned marker column (col H) lined up against
# every row of the (now taller) Manchester block.
$ws.Range("H4:H7").HorizontalAlignment = -4152

# Re-create the two mailto hyperlinks at their new row positions.
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Manchesteret@justice.gov.uk", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "Manchesteret@justice.gov.uk") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:glasgowet@justice.gov.uk", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "glasgowet@justice.gov.uk") | Out-Null

# Match the author's final cursor position.
$ws.Range("A6").Select() | Out-Null
